$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "addind preparer to sheet" - update the purpose column (E2:E7) from
# "S.GISH" to the newly introduced "fullRNASEQ" value.
$ws.Range("E2:E7").Value = "fullRNASEQ"

# Match the cursor/selection left behind by the editing session.
[void]$ws.Range("D8:F13").Select()
